$wb = $excel.ActiveWorkbook

# --- Editor sheet: add a new "expectedOutputType" column (C) with data ---
$editor = $wb.Worksheets.Item("Editor")

$editor.Range("C1").Value = "expectedOutputType"
$editor.Range("C2").Value = "Alert"
$editor.Range("C3").Value = "Console"
$editor.Range("C4").Value = "Alert"
$editor.Range("C5").Value = "Alert"

# Target OOXML column width is 18.7265625 characters; the host's
# ColumnWidth setter only resolves to pixel-quantised values (multiples of
# 1/6), so 17.8 is chosen as the input that lands on the closest
# achievable width to the target.
$editor.Columns.Item(3).ColumnWidth = 17.8

# Make Editor the active/selected sheet with the new active cell
$editor.Activate()
$editor.Range("C8").Select()

# --- Register sheet: it was previously the tab-selected sheet; ---
# activating Editor above already moves tab selection away from Register.
